{"js": "// Replace the date line and each two-digit multiplication problem with\n// its new value. Every \"before\" text below is unique within the\n// document, so a simple global search/replace for each pair is safe.\nconst replacements = [\n  [\"2025-08-17 Sunday\", \"2025-08-18 Monday\"],\n  [\"59\u00d736=\", \"31\u00d790=\"],\n  [\"60\u00d731=\", \"73\u00d713=\"],\n  [\"46\u00d763=\", \"71\u00d726=\"],\n  [\"16\u00d713=\", \"83\u00d780=\"],\n  [\"56\u00d762=\", \"79\u00d714=\"],\n  [\"48\u00d785=\", \"19\u00d771=\"],\n  [\"31\u00d730=\", \"88\u00d797=\"],\n  [\"84\u00d776=\", \"12\u00d786=\"],\n  [\"16\u00d738=\", \"64\u00d787=\"],\n  [\"47\u00d767=\", \"14\u00d749=\"],\n  [\"27\u00d767=\", \"81\u00d730=\"],\n  [\"75\u00d731=\", \"56\u00d736=\"],\n  [\"47\u00d756=\", \"11\u00d732=\"],\n  [\"26\u00d789=\", \"90\u00d776=\"],\n  [\"87\u00d777=\", \"68\u00d747=\"],\n  [\"21\u00d725=\", \"70\u00d733=\"],\n  [\"63\u00d729=\", \"26\u00d782=\"],\n  [\"48\u00d760=\", \"45\u00d783=\"],\n  [\"97\u00d711=\", \"45\u00d721=\"],\n  [\"92\u00d717=\", \"25\u00d792=\"],\n  [\"38\u00d762=\", \"60\u00d730=\"],\n  [\"62\u00d780=\", \"32\u00d771=\"],\n  [\"55\u00d739=\", \"81\u00d759=\"],\n  [\"38\u00d731=\", \"34\u00d716=\"],\n  [\"36\u00d758=\", \"26\u00d717=\"],\n];\n\nconst body = context.document.body;\n\nfor (const [before, after] of replacements) {\n  const results = body.search(before, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (const range of results.items) {\n    range.insertText(after, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Replace the date line and each two-digit multiplication problem with\n# its new value. Every \"before\" text below is unique within the\n# document, so a simple global Find/Replace for each pair is safe.\n$d = $word.ActiveDocument\n\n$pairs = @(\n  @(\"2025-08-17 Sunday\", \"2025-08-18 Monday\"),\n  @(\"59\u00d736=\", \"31\u00d790=\"),\n  @(\"60\u00d731=\", \"73\u00d713=\"),\n  @(\"46\u00d763=\", \"71\u00d726=\"),\n  @(\"16\u00d713=\", \"83\u00d780=\"),\n  @(\"56\u00d762=\", \"79\u00d714=\"),\n  @(\"48\u00d785=\", \"19\u00d771=\"),\n  @(\"31\u00d730=\", \"88\u00d797=\"),\n  @(\"84\u00d776=\", \"12\u00d786=\"),\n  @(\"16\u00d738=\", \"64\u00d787=\"),\n  @(\"47\u00d767=\", \"14\u00d749=\"),\n  @(\"27\u00d767=\", \"81\u00d730=\"),\n  @(\"75\u00d731=\", \"56\u00d736=\"),\n  @(\"47\u00d756=\", \"11\u00d732=\"),\n  @(\"26\u00d789=\", \"90\u00d776=\"),\n  @(\"87\u00d777=\", \"68\u00d747=\"),\n  @(\"21\u00d725=\", \"70\u00d733=\"),\n  @(\"63\u00d729=\", \"26\u00d782=\"),\n  @(\"48\u00d760=\", \"45\u00d783=\"),\n  @(\"97\u00d711=\", \"45\u00d721=\"),\n  @(\"92\u00d717=\", \"25\u00d792=\"),\n  @(\"38\u00d762=\", \"60\u00d730=\"),\n  @(\"62\u00d780=\", \"32\u00d771=\"),\n  @(\"55\u00d739=\", \"81\u00d759=\"),\n  @(\"38\u00d731=\", \"34\u00d716=\"),\n  @(\"36\u00d758=\", \"26\u00d717=\")\n)\n\nforeach ($pair in $pairs) {\n  $before = $pair[0]\n  $after = $pair[1]\n\n  $find = $d.Content.Find\n  $find.ClearFormatting()\n  $find.Replacement.ClearFormatting()\n  $find.Text = $before\n  $find.Replacement.Text = $after\n  $find.Execute([ref]$before, $false, $false, $false, $false, $false, $true, 0, $false, [ref]$after, 2) | Out-Null\n}\n"}
